# This edit performs a 3-row cyclic rotation of the observation records
# currently sitting in rows 2, 3 and 4 of the active sheet:
#   new row 2  <-  old row 4
#   new row 3  <-  old row 2
#   new row 4  <-  old row 3
#
# Only the columns whose content actually differs between the three rows
# are touched (A, B, E, F, G, H, J, L, M, Q, R, S, AC, AF) - the remaining
# columns hold identical data in all three rows, so rotating them would be
# a no-op (and risks corrupting date-looking text cells via Excel's value
# auto-coercion), so we deliberately leave them alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (1-based) that differ between rows 2/3/4.
$cols = @(1, 2, 5, 6, 7, 8, 10, 12, 13, 17, 18, 19, 29, 32)

$row2 = @{}
$row3 = @{}
$row4 = @{}

foreach ($c in $cols) {
    $row2[$c] = $ws.Cells.Item(2, $c).Value()
    $row3[$c] = $ws.Cells.Item(3, $c).Value()
    $row4[$c] = $ws.Cells.Item(4, $c).Value()
}

# Write the rotated values into place.
foreach ($c in $cols) {
    $ws.Cells.Item(2, $c).Value = $row4[$c]
    $ws.Cells.Item(3, $c).Value = $row2[$c]
    $ws.Cells.Item(4, $c).Value = $row3[$c]
}

# Columns 10 (J) and 12 (L) are blank placeholder cells (present, but with
# no content) that need to exist in specific destination rows after the
# rotation. Setting .Value = "" clears/removes the cell outright instead of
# leaving an empty text cell behind, so use the classic leading-apostrophe
# trick to force an empty *text* cell, then reset the style Excel applies
# for the quote-prefix back to Normal so no stray formatting is left.
$blankCells = @("L2", "J3", "AF3", "L4")
foreach ($addr in $blankCells) {
    $ws.Range($addr).Value = "'"
    $ws.Range($addr).Style = "Normal"
}

Write-Output "Rotated rows 2-4 successfully"
